$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" -------------
# Overview sheet: columns E (zh-cn status) and F (de-de status), rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column C, rows 2-3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# de-de sheet: Status column C, rows 2-3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the status columns to fit the shorter text ----------------------
# Target stored width is ~13.41 (previously ~17.22); the COM ColumnWidth
# setter here snaps to 1/6-character increments, so feed it the input that
# lands on the closest achievable value to the recorded width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
